# Applies the "data up otto 14th" commit: extends the daily survey table
# (Sheet1) with seven more date-rows (08-14 Sep 2020) and fixes a small
# batch of already-entered values in the two days before that (rows 217-220).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: small corrections to existing cells (rows 217-220) ---
# Each entry is row, column, corrected value.
$tweaks = @(
    @(217, 12, 0.6431656),   # L217
    @(217, 24, 0.339749),   # X217
    @(217, 55, 0.5161563),   # BC217
    @(218, 12, 0.6206581),   # L218
    @(218, 24, 0.3407038),   # X218
    @(218, 55, 0.5075484),   # BC218
    @(219, 12, 0.6247034),   # L219
    @(219, 24, 0.3445838),   # X219
    @(219, 55, 0.5466399),   # BC219
    @(220, 7, 0.4528665),   # G220
    @(220, 12, 0.594932),   # L220
    @(220, 24, 0.3162993),   # X220
    @(220, 49, 0.7565516),   # AW220
    @(220, 55, 0.5265579)   # BC220
)

foreach ($t in $tweaks) {
    $ws.Cells.Item($t[0], $t[1]).Value = $t[2]
}

# --- Part 2: seven new daily rows, 221-228 ---
# Row 221 (07 09 2020) was already present with just its date label; it now
# gets its per-state figures. Rows 222-227 (08-13 Sep 2020) are brand new,
# fully populated rows. Row 228 (14 09 2020) is the newest day and - like
# every other "just added" row in this sheet - starts out with only its
# date label in column A and no data yet.
#
# Columns E, N, AC, AR and AZ are intentionally left blank in every row of
# this sheet (those territories are not surveyed), so they are skipped below.

$newRowData = @(
    @{ Row = 221; Date = "07 09 2020"; Values = @{ 2=0.3347136; 3=1.0544599; 4=0.8910583; 6=0.5968199; 7=0.4252703; 8=0.5148606999999999; 9=0.2166225; 10=0.6474417; 11=0.3277264; 12=0.6197112; 13=0.8001451000000001; 15=0.1966146; 16=0.991238; 17=0.9885546; 18=0.6243072; 19=0.7867527; 20=0.7422598; 21=0.6093631; 22=0.6904088; 23=0.3211003; 24=0.2889477; 25=0.3503425; 26=0.4554259; 27=0.6122507; 28=0.8850158; 30=1.0218092; 31=0.4746033; 32=0.5224618; 33=0.7389217; 34=1.3493849; 35=0.2810693; 36=0.1867376; 37=0.4197149; 38=0.76176; 39=0.3840781; 40=0.4520036; 41=0.9354709; 42=0.3392905; 43=0.3720495; 45=0.598445; 46=0.6198005; 47=1.2248118; 48=0.8262567; 49=0.7203092; 50=0.7059873; 51=0.5209808; 53=0.621519; 54=0.4768944; 55=0.5624344999999999; 56=0.516447; 57=0.9531287000000001 } },
    @{ Row = 222; Date = "08 09 2020"; Values = @{ 2=0.4226986; 3=0.94617; 4=0.7737508; 6=0.5894885; 7=0.4118265; 8=0.5607722000000001; 9=0.1904293; 10=0.4937092; 11=0.3061224; 12=0.6270478; 13=0.8459105; 15=0.2010117; 16=1.0553306; 17=0.9245186; 18=0.6633119; 19=0.8245114; 20=0.7105188; 21=0.5195399000000001; 22=0.7951041; 23=0.3124991; 24=0.324725; 25=0.3294777; 26=0.4740468; 27=0.5820455; 28=0.9334741; 30=0.9548885; 31=0.5393911; 32=0.5597391; 33=0.8451945; 34=1.3441838; 35=0.2075318; 36=0.1795839; 37=0.3902468; 38=0.7643867; 39=0.3671089; 40=0.5071397; 41=0.9462206; 42=0.3680947; 43=0.3928694; 45=0.6525639; 46=0.679525; 47=1.2697386; 48=0.7489366; 49=0.7061197; 50=0.6680397; 51=0.5931243; 53=0.4992088; 54=0.4657631; 55=0.5743646; 56=0.6682018; 57=0.8707214 } },
    @{ Row = 223; Date = "09 09 2020"; Values = @{ 2=0.4694646; 3=0.8980487; 4=0.8660442; 6=0.5562953; 7=0.3983495; 8=0.6343137; 9=0.1791058; 10=0.1110289; 11=0.4022072; 12=0.6197759; 13=0.7849433; 15=0.2982601; 16=1.0133926; 17=0.9828504; 18=0.6312591; 19=0.8053696; 20=0.7018316; 21=0.5152156; 22=0.7404828; 23=0.3522254; 24=0.2497851; 25=0.2711324; 26=0.4340591; 27=0.5389932; 28=0.9918337; 30=0.8786075; 31=0.5427902999999999; 32=0.5914454; 33=0.7340644; 34=1.0387021; 35=0.2034429; 36=0.2290289; 37=0.3194227; 38=0.6148432; 39=0.3969299; 40=0.4768183; 41=0.9983642; 42=0.2887357; 43=0.3249615; 45=0.6289955; 46=0.6528433; 47=1.3152443; 48=0.675307; 49=0.6907435; 50=0.6043153999999999; 51=0.5994249; 53=0.5477065; 54=0.4101413; 55=0.4824523; 56=0.6406212; 57=1.048486 } },
    @{ Row = 224; Date = "10 09 2020"; Values = @{ 2=0.4246616; 3=0.8847853; 4=0.8026817000000001; 6=0.5531236; 7=0.3884633; 8=0.5804665; 9=0.1681555; 10=0; 11=0.4952713; 12=0.6031571; 13=0.7616379; 15=0.3048264; 16=1.0478427; 17=0.8946232; 18=0.6367683; 19=0.7716912; 20=0.7111976; 21=0.4972781; 22=0.8185757; 23=0.2972767; 24=0.2811277; 25=0.3592979; 26=0.4684673; 27=0.5375185; 28=0.9580867; 30=1.0232772; 31=0.6165333; 32=0.5553960999999999; 33=0.6913122; 34=0.9635994; 35=0.1518212; 36=0.2227696; 37=0.2653212; 38=0.7148747; 39=0.3800577; 40=0.5082738999999999; 41=0.9332014; 42=0.3820252; 43=0.3252678; 45=0.6472879; 46=0.6646465; 47=1.1500744; 48=0.7169183; 49=0.7275092; 50=0.6311781; 51=0.5121582; 53=0.5423761; 54=0.4171337; 55=0.4188867; 56=0.5580549; 57=0.940327 } },
    @{ Row = 225; Date = "11 09 2020"; Values = @{ 2=0.3784219; 3=0.8018313; 4=0.760732; 6=0.4868956; 7=0.3811755; 8=0.6331548; 9=0.1909445; 10=0; 11=0.4284263; 12=0.581942; 13=0.7176731; 15=0.3978526; 16=1.003167; 17=0.7569630000000001; 18=0.6074746; 19=0.698183; 20=0.7535865; 21=0.4753264; 22=0.8046473; 23=0.342143; 24=0.3268361; 25=0.3068984; 26=0.4255997; 27=0.5718358; 28=0.8951195; 30=1.0176151; 31=0.6812729; 32=0.547832; 33=0.8024391; 34=0.7890547; 35=0.1652841; 36=0.2462155; 37=0.2065955; 38=0.7232819; 39=0.3742921; 40=0.5553046; 41=0.8974084; 42=0.3527734; 43=0.3073712; 45=0.6628324; 46=0.630559; 47=1.125204; 48=0.7695033; 49=0.7337985; 50=0.5238111; 51=0.5030476; 53=0.4706698; 54=0.4137233; 55=0.4498924; 56=0.6023816; 57=0.5462518 } },
    @{ Row = 226; Date = "12 09 2020"; Values = @{ 2=0.376438; 3=0.6686873; 4=0.7058416; 6=0.484007; 7=0.3886581; 8=0.6304933; 9=0.1702413; 10=0; 11=0.4287773; 12=0.6180247; 13=0.6635294; 15=0.3198504; 16=1.0034107; 17=0.7229949; 18=0.544471; 19=0.6886544999999999; 20=0.7302925; 21=0.5167663; 22=0.648145; 23=0.340676; 24=0.3243584; 25=0.3067318; 26=0.4101403; 27=0.5643594; 28=0.9390447; 30=0.9266152; 31=0.5073671; 32=0.5228221; 33=0.7867774; 34=0.5953083; 35=0.0850403; 36=0.2484598; 37=0.2074978; 38=0.7231595; 39=0.3434071; 40=0.5626038; 41=0.8265842; 42=0.3233294; 43=0.3835625; 45=0.7312332; 46=0.6831153; 47=1.1214745; 48=0.7458314; 49=0.6531695; 50=0.5185835; 51=0.4600205; 53=0.4763337; 54=0.4242207; 55=0.4950351; 56=0.5613345; 57=0.7340801 } },
    @{ Row = 227; Date = "13 09 2020"; Values = @{ 2=0.4468255; 3=0.758768; 4=0.8387426; 6=0.4729999; 7=0.3669905; 8=0.5455095; 9=0.2173014; 10=0; 11=0.4064594; 12=0.6177002; 13=0.6227791; 15=0.3261091; 16=0.9180692; 17=0.7216827; 18=0.5213421; 19=0.680103; 20=0.64749; 21=0.433487; 22=0.6660688; 23=0.2884764; 24=0.3418484; 25=0.3663912; 26=0.3886767; 27=0.591836; 28=0.9505425; 30=0.9110662; 31=0.5333203; 32=0.5791752; 33=0.7361115; 34=0.6764513; 35=0.1178698; 36=0.2540021; 37=0.1238624; 38=0.6716746; 39=0.2949663; 40=0.5975058; 41=0.8903515; 42=0.3900878; 43=0.381106; 45=0.4909864; 46=0.6824189000000001; 47=1.085925; 48=0.6867671; 49=0.5884046000000001; 50=0.560606; 51=0.4229261; 53=0.3356003; 54=0.3794238; 55=0.5349484; 56=0.5859354; 57=0.8692054 } },
    @{ Row = 228; Date = "14 09 2020"; Values = @{} }
)

foreach ($rowEntry in $newRowData) {
    $ws.Cells.Item($rowEntry.Row, 1).Value = $rowEntry.Date
    foreach ($colNum in $rowEntry.Values.Keys) {
        $ws.Cells.Item($rowEntry.Row, [int]$colNum).Value = $rowEntry.Values[$colNum]
    }
}
